# "update excel & package name fix"
# Fill in the POC (point of contact) names that were missing for two rows
# of the task-distribution table on Sheet1, and move the active-cell
# selection down one row (B3 -> B4), matching the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("D3").Value = "Ken"
$ws.Range("D11").Value = "Sunny"

$ws.Range("B4").Select()
